$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append starting at row 14 (Country, Problem described)
$rows = @(
    @{ A = "Eritrea";        B = "Issues with getting the data from OSM. Cables again" },
    @{ A = "Ghana";          B = "Issue with Nan numbers when building shapes. Also the OSM data seems to span multiple countries." },
    @{ A = "Togo";           B = "OSM issue with cables" },
    @{ A = "Burkina Fasso";  B = "OSM issue with cables" },
    @{ A = "Sierra Leone";   B = "OSM issue with cables" },
    @{ A = "Guinea";         B = "OSM issue with cables" },
    @{ A = "Guinea Bissau";  B = "OSM issue with cables" },
    @{ A = "Western Sahara"; B = "No powerplants so ppmatching breaks" },
    @{ A = "Madgascar";      B = "OSM issue with cables" },
    @{ A = "Malawi ";        B = "OSM issue with cables" }
)

$heights = @(48, 96, 16, 16, 16, 16, 16, 32, 16, 16)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i].A
    $ws.Cells.Item($r, 2).Value = $rows[$i].B
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 2)).WrapText = $true
    $ws.Rows.Item($r).RowHeight = $heights[$i]
}

# Update selection/view to reflect the new active cell (B23) and scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("B23").Select()
